# "Included option for ECS to have a dumb-bob prism hanging down below"
#
# 1) Insert a new worksheet "Mars Eclipse" right after "Mars-Mars" and
#    populate it with the Mars-diameter / orbit-circumference / eclipse
#    half-angle (theta) table.
# 2) On "Mars-Mars", change the two orbital-altitude inputs (C7, C13) from
#    12000 to 7000 km; every other changed cell on that sheet is a formula
#    that recalculates automatically.
# 3) Restore view/selection state to match: Mars-Mars keeps the active tab
#    with selection E26; the new sheet's own selection is M11.

$wb = $excel.ActiveWorkbook

$marsMars = $wb.Worksheets.Item("Mars-Mars")

# --- 1. New "Mars Eclipse" sheet, inserted right after "Mars-Mars" -------
$eclipse = $wb.Worksheets.Add($null, $marsMars)
$eclipse.Name = "Mars Eclipse"

# Row 4 headers and the Mars-diameter label/value are entered in the same
# order the author typed them so the shared-string table comes out in the
# same order: "Orbit SMA", "Circumference", "Mars diam", "Theta".
$eclipse.Range("C4").Value = "Orbit SMA"
$eclipse.Range("D4").Value = "Circumference"
$eclipse.Range("C2").Value = "Mars diam"
$eclipse.Range("E4").Value = "Theta"

$eclipse.Range("D2").Formula = '=2*3390'

# Data rows 5-36: orbit SMA 5000..20500 (step 500), circumference,
# eclipse half-angle theta, and the fraction of the orbit spent eclipsed.
for ($i = 0; $i -lt 32; $i++) {
    $row = 5 + $i
    $sma = 5000 + $i * 500
    $eclipse.Cells.Item($row, 3).Value = $sma
    $eclipse.Cells.Item($row, 4).Formula = "=2*C$row*PI()"
    $eclipse.Cells.Item($row, 4).NumberFormat = "0"
    $eclipse.Cells.Item($row, 5).Formula = '=2*ASIN($D$2/(2*C' + $row + '))'
    $eclipse.Cells.Item($row, 5).NumberFormat = "0.00"
    $eclipse.Cells.Item($row, 6).Formula = "=E$row/(2*PI())"
}

# --- 2. "Mars-Mars": drop the two periapsis-altitude inputs to 7000 km ---
$marsMars.Range("C7").Value = 7000
$marsMars.Range("C13").Value = 7000

# --- 3. View/selection state ---------------------------------------------
$eclipse.Range("M11").Select()
$marsMars.Select()
$marsMars.Range("E26").Select()
